$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Rebuild the whole text body in one shot (paragraphs separated by CR).
# This inherits the original paragraph-1 run formatting (lang="fr-CA", sz=2400)
# for every new paragraph/run, which we then touch up below.
$tr.Text = "Faire les analyses quantitatives ayant mene a la publication suivante`rR NG_Lab - Tutoriel et resultat SAE3D.pdf`rLes instructions :`rsae3d_tutorial.html`rLes donnees (simulees!) :`rsae3d.csv`rLe code en exemple :`rDiapos subsequentes"

# Fix the accented characters that are awkward to embed directly in the script.
$tr.Paragraphs(1,1).Text = "Faire les analyses quantitatives ayant mené à la publication suivante"
$tr.Paragraphs(5,1).Text = "Les données (simulées!) :"

# Paragraph-level indent (lvl="1" == IndentLevel 2) for the "sub" bullet paragraphs.
$tr.Paragraphs(2,1).IndentLevel = 2
$tr.Paragraphs(4,1).IndentLevel = 2
$tr.Paragraphs(6,1).IndentLevel = 2
$tr.Paragraphs(8,1).IndentLevel = 2

# Paragraph 2: "R NG_Lab - Tutoriel et résultat SAE3D.pdf" -> sz=2000, italic, lang fr-FR
$p2 = $tr.Paragraphs(2,1)
$p2.Text = "R NG_Lab - Tutoriel et résultat SAE3D.pdf"
$p2.Font.Size = 20
$p2.Font.Italic = 1

# Paragraph 4: "sae3d_tutorial.html" -> sz=2000, italic
$p4 = $tr.Paragraphs(4,1)
$p4.Font.Size = 20
$p4.Font.Italic = 1

# Paragraph 6: "sae3d.csv" -> sz=2000, italic
$p6 = $tr.Paragraphs(6,1)
$p6.Font.Size = 20
$p6.Font.Italic = 1

# Paragraph 8: "Diapos subséquentes" -> sz=2000, italic
$p8 = $tr.Paragraphs(8,1)
$p8.Text = "Diapos subséquentes"
$p8.Font.Size = 20
$p8.Font.Italic = 1

# Picture 5 (id=6): nudge its position.
$pic = $s.Shapes.Item(4)
$pic.Left = 427.2227783203125
$pic.Top = 143.75
